$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rangeAddr, [string]$text)
    $escapedForFormula = $text -replace '"', '""'
    $ws.Range($rangeAddr).Formula = '="' + $escapedForFormula + '"'
    $ws.Range($rangeAddr).Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)
}

Set-TextValue 'D2' '23.547.16'
Set-TextValue 'E2' '  +1.50%  '
Set-TextValue 'D3' '1.640.28'
Set-TextValue 'E3' '  +2.46%  '
Set-TextValue 'D4' '0.9998'
Set-TextValue 'E4' '  +0.25%  '
Set-TextValue 'D5' '308.05'
Set-TextValue 'E5' '  +1.81%  '
Set-TextValue 'D6' '0.9995'
Set-TextValue 'E6' '  +0.13%  '
Set-TextValue 'D7' '0.3758'
Set-TextValue 'E7' '  -0.50%  '
Set-TextValue 'D8' '52.73'
Set-TextValue 'E8' '  +3.10%  '
Set-TextValue 'D9' '0.3671'
Set-TextValue 'E9' '  +1.84%  '
Set-TextValue 'D10' '1.279'
Set-TextValue 'E10' '  +1.42%  '
Set-TextValue 'E11' '  +1.02%  '
Set-TextValue 'D12' '0.9998'
Set-TextValue 'E12' '  +0.27%  '
Set-TextValue 'D13' '23.04'
Set-TextValue 'E13' '  +1.78%  '
Set-TextValue 'D14' '6.676'
Set-TextValue 'E14' '  +1.50%  '
Set-TextValue 'D15' '0.00001285'
Set-TextValue 'E15' '  +2.83%  '
Set-TextValue 'D16' '7.458'
Set-TextValue 'E16' '  +1.18%  '
Set-TextValue 'D17' '1.641.99'
Set-TextValue 'E17' '  +2.64%  '
Set-TextValue 'D18' '95.11'
Set-TextValue 'E18' '  +1.74%  '
Set-TextValue 'D19' '0.06919'
Set-TextValue 'E19' '  +1.07%  '
Set-TextValue 'D20' '18.30'
Set-TextValue 'E20' '  +1.59%  '
Set-TextValue 'D21' '6.581'
Set-TextValue 'E21' '  +0.90%  '
Set-TextValue 'D22' '0.9984'
Set-TextValue 'E22' '  -0.09%  '
Set-TextValue 'D23' '23.541.55'
Set-TextValue 'E23' '  +1.47%  '
Set-TextValue 'D24' '12.90'
Set-TextValue 'E24' '  -0.14%  '
Set-TextValue 'D25' '3.086'
Set-TextValue 'E25' '  +3.79%  '
Set-TextValue 'D26' '2.420'
Set-TextValue 'E26' '  +0.96%  '
Set-TextValue 'D27' '21.33'
Set-TextValue 'E27' '  +0.94%  '
Set-TextValue 'D28' '151.94'
Set-TextValue 'E28' '  +1.77%  '
Set-TextValue 'D29' '5.345'
Set-TextValue 'E29' '  +2.07%  '
Set-TextValue 'D30' '136.43'
Set-TextValue 'E30' '  +2.13%  '
Set-TextValue 'D31' '2.376'
Set-TextValue 'E31' '  -0.54%  '
Set-TextValue 'D32' '1.827.03'
Set-TextValue 'E32' '  +3.08%  '
Set-TextValue 'D33' '6.875'
Set-TextValue 'E33' '  +0.50%  '
Set-TextValue 'D34' '0.9811'
Set-TextValue 'E34' '  +0.14%  '
Set-TextValue 'D35' '0.02857'
Set-TextValue 'E35' '  +5.79%  '
Set-TextValue 'D36' '10.45'
Set-TextValue 'E36' '  +1.52%  '
Set-TextValue 'D37' '0.07421'
Set-TextValue 'E37' '  -1.77%  '
Set-TextValue 'D38' '0.2561'
Set-TextValue 'E38' '  +2.48%  '
Set-TextValue 'D39' '6.213'
Set-TextValue 'E39' '  +0.83%  '
Set-TextValue 'D40' '0.08904'
Set-TextValue 'E40' '  +1.27%  '
Set-TextValue 'D41' '1.386'
Set-TextValue 'E41' '  +1.60%  '
Set-TextValue 'D42' '0.7138'
Set-TextValue 'E42' '  +0.24%  '
Set-TextValue 'B43' 'EnergySwap'
Set-TextValue 'C43' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D43' '16.37'
Set-TextValue 'E43' '  +5.58%  '
Set-TextValue 'B44' 'Aptos'
Set-TextValue 'C44' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D44' '12.58'
Set-TextValue 'E44' '  +1.36%  '
Set-TextValue 'D45' '0.6581'
Set-TextValue 'E45' '  +1.08%  '
Set-TextValue 'D46' '2.355'
Set-TextValue 'E46' '  +2.37%  '
Set-TextValue 'D47' '4.049'
Set-TextValue 'E47' '  +1.00%  '
Set-TextValue 'D48' '0.9987'
Set-TextValue 'E48' '  +0.13%  '
Set-TextValue 'D49' '0.08002'
Set-TextValue 'E49' '  +0.85%  '
Set-TextValue 'D50' '130.20'
Set-TextValue 'E50' '  -1.23%  '
Set-TextValue 'E51' '  +0.45%  '

$excel.CutCopyMode = 0
Write-Output "Done applying updates."
